$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue "D2" '64.957.97'
Set-TextValue "E2" '  +0.45%  '

# Row 3
Set-TextValue "D3" '3.524.80'
Set-TextValue "E3" '  +3.36%  '

# Row 4
Set-TextValue "E4" '  +0.46%  '

# Row 5
Set-TextValue "D5" '594.98'
Set-TextValue "E5" '  +2.05%  '

# Row 6
Set-TextValue "D6" '136.50'
Set-TextValue "E6" '  +0.42%  '

# Row 7
Set-TextValue "D7" '3.524.11'
Set-TextValue "E7" '  +3.19%  '

# Row 8
Set-TextValue "E8" '  +0.31%  '

# Row 9
Set-TextValue "E9" '  +1.14%  '

# Row 10
Set-TextValue "E10" '  +2.19%  '

# Row 11
Set-TextValue "E11" '  -1.74%  '

# Row 12
Set-TextValue "D12" '0.383'
Set-TextValue "E12" '  +2.61%  '

# Row 13
Set-TextValue "D13" '4.129.04'
Set-TextValue "E13" '  +3.77%  '

# Row 14
Set-TextValue "B14" 'Avalanche'
Set-TextValue "C14" 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue "D14" '27.10'
Set-TextValue "E14" '  +3.42%  '

# Row 15
Set-TextValue "B15" 'ShibaInu'
Set-TextValue "C15" 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue "D15" '0.0000180'
Set-TextValue "E15" '  +2.39%  '

# Row 16
Set-TextValue "D16" '3.528.61'
Set-TextValue "E16" '  +3.54%  '

# Row 17
Set-TextValue "E17" '  +1.41%  '

# Row 18
Set-TextValue "D18" '64.961.22'
Set-TextValue "E18" '  +0.78%  '

# Row 19
Set-TextValue "D19" '10.00'
Set-TextValue "E19" '  +4.86%  '

# Row 20
Set-TextValue "D20" '5.82'
Set-TextValue "E20" '  +0.45%  '

# Row 21
Set-TextValue "D21" '14.14'
Set-TextValue "E21" '  +4.98%  '

# Row 22
Set-TextValue "D22" '387.87'
Set-TextValue "E22" '  +1.71%  '

# Row 23
Set-TextValue "D23" '0.572'
Set-TextValue "E23" '  +4.27%  '

# Row 24
Set-TextValue "D24" '3.668.06'
Set-TextValue "E24" '  +3.67%  '

# Row 25
Set-TextValue "D25" '73.80'
Set-TextValue "E25" '  +2.83%  '

# Row 26
Set-TextValue "D26" '0.999'

# Row 27
Set-TextValue "D27" '0.0000112'
Set-TextValue "E27" '  +8.60%  '

# Row 28
Set-TextValue "D28" '7.63'
Set-TextValue "E28" '  +7.13%  '

# Row 29
Set-TextValue "D29" '0.999'
Set-TextValue "E29" '  -0.42%  '

# Row 30
Set-TextValue "E30" '  +3.33%  '

# Row 31
Set-TextValue "D31" '8.15'
Set-TextValue "E31" '  +1.69%  '

# Row 32
Set-TextValue "D32" '3.542.19'
Set-TextValue "E32" '  +3.69%  '

# Row 33
Set-TextValue "E33" '  +0.08%  '

# Row 34
Set-TextValue "D34" '23.67'
Set-TextValue "E34" '  +3.27%  '

# Row 35
Set-TextValue "D35" '1.34'
Set-TextValue "E35" '  +13.53%  '

# Row 36
Set-TextValue "E36" '  +3.15%  '

# Row 37
Set-TextValue "D37" '169.84'
Set-TextValue "E37" '  +1.27%  '

# Row 38
Set-TextValue "E38" '  +6.67%  '

# Row 39
Set-TextValue "D39" '6.80'
Set-TextValue "E39" '  +1.27%  '

# Row 40
Set-TextValue "D40" '4.93'
Set-TextValue "E40" '  +7.40%  '

# Row 41
Set-TextValue "D41" '0.0795'
Set-TextValue "E41" '  +5.75%  '

# Row 42
Set-TextValue "E42" '  +1.48%  '

# Row 43
Set-TextValue "D43" '26.43'
Set-TextValue "E43" '  +17.40%  '

# Row 44
Set-TextValue "B44" 'FirstDigitalUSD'
Set-TextValue "C44" 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue "D44" '1.00'
Set-TextValue "E44" '  +0.50%  '

# Row 45
Set-TextValue "B45" 'OKB'
Set-TextValue "C45" 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue "D45" '42.44'
Set-TextValue "E45" '  +0.59%  '

# Row 46
Set-TextValue "D46" '4.40'
Set-TextValue "E46" '  +2.67%  '

# Row 47
Set-TextValue "D47" '1.19'
Set-TextValue "E47" '  +6.32%  '

# Row 48
Set-TextValue "D48" '1.65'
Set-TextValue "E48" '  +4.50%  '

# Row 49
Set-TextValue "E49" '  +6.44%  '

# Row 50
Set-TextValue "D50" '2.403.52'
Set-TextValue "E50" '  +11.23%  '

# Row 51
Set-TextValue "D51" '301.07'
Set-TextValue "E51" '  +9.70%  '
